# Update the "想去人数" (number of people interested) figures for the
# first three con listings on both the "展览" and "全部类型" sheets,
# matching the refreshed data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 730
    $ws.Range("F3").Value = 4104
    $ws.Range("F4").Value = 116
}
